{"js": "// The document contains a \"Diario de Reflexi\u00f3n\" table. In the answer to the\n// first question (\"\u00bfCu\u00e1les son las asignaturas o certificados que m\u00e1s te\n// gustaron...?\") the three paragraphs that hold the author's reflection text\n// were left with default (left) alignment, while the rest of the document's\n// answer paragraphs are already justified. This change brings those three\n// paragraphs in line by setting their alignment to \"justify\" (OOXML\n// <w:jc w:val=\"both\"/>).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Unique, stable snippets of text from each of the three target paragraphs.\nconst targetSnippets = [\n  \"Las asignaturas iniciales de modelamiento, consulta y programaci\u00f3n de base de datos (Oracle) me gustaron mucho\",\n  \"Las asignaturas de programaci\u00f3n web, m\u00f3vil, aplicaciones de escritorio con Java y arquitectura fueron muy motivadoras\",\n  \"Mientras que las asignaturas de ML me interesaron y dieron una visi\u00f3n diferente\"\n];\n\nlet found = 0;\nfor (const p of paragraphs.items) {\n  const text = p.text || \"\";\n  if (targetSnippets.some((snippet) => text.indexOf(snippet) !== -1)) {\n    p.alignment = Word.Alignment.justified;\n    found++;\n  }\n}\n\nawait context.sync();\n\nif (found !== targetSnippets.length) {\n  throw new Error(\n    \"Expected to justify \" + targetSnippets.length + \" paragraphs, found \" + found\n  );\n}\n", "ps1": "# The document contains a \"Diario de Reflexi\u00f3n\" table. In the answer to the\n# first question (\"\u00bfCu\u00e1les son las asignaturas o certificados que m\u00e1s te\n# gustaron...?\") the three paragraphs that hold the author's reflection text\n# were left with default (left) alignment, while the rest of the document's\n# answer paragraphs are already justified. This change brings those three\n# paragraphs in line by setting their alignment to justify\n# (wdAlignParagraphJustify = 3, i.e. OOXML <w:jc w:val=\"both\"/>).\n\n$d = $word.ActiveDocument\n\n# Unique, stable snippets of text from each of the three target paragraphs.\n$targets = @(\n    \"Las asignaturas iniciales de modelamiento, consulta y programaci\u00f3n de base de datos (Oracle) me gustaron mucho\",\n    \"Las asignaturas de programaci\u00f3n web, m\u00f3vil, aplicaciones de escritorio con Java y arquitectura fueron muy motivadoras\",\n    \"Mientras que las asignaturas de ML me interesaron y dieron una visi\u00f3n diferente\"\n)\n\n$wdAlignParagraphJustify = 3\n$changed = 0\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    foreach ($snippet in $targets) {\n        if ($text -like \"*$snippet*\") {\n            $p.Alignment = $wdAlignParagraphJustify\n            $changed++\n            break\n        }\n    }\n}\n\nif ($changed -ne $targets.Count) {\n    throw \"Expected to justify $($targets.Count) paragraphs, justified $changed\"\n}\n\nWrite-Output \"Justified $changed paragraph(s)\"\n"}
